$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "260.70"

# Row 4
Set-TextValue 4 4 "6.211"

# Row 5
Set-TextValue 5 4 "0.06075"

# Row 7
Set-TextValue 7 4 "6.711"

# Row 8
Set-TextValue 8 4 "1.358"

# Row 9
Set-TextValue 9 4 "0.7991"

# Row 11
Set-TextValue 11 4 "0.08126"

# Row 12
Set-TextValue 12 4 "0.03311"

# Row 13
Set-TextValue 13 4 "0.03122"

# Row 14
Set-TextValue 14 4 "0.09267"

# Row 15
Set-TextValue 15 4 "3.895"

# Row 16
Set-TextValue 16 4 "0.001691"

# Row 17
Set-TextValue 17 4 "0.04822"

# Row 18
$ws.Cells.Item(18, 2).Value = "TigerCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue 18 4 "0.006191"
$ws.Cells.Item(18, 5).Value = "17TigerCashTCH"

# Row 19
$ws.Cells.Item(19, 2).Value = "BitKan"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue 19 4 "0.001100"
$ws.Cells.Item(19, 5).Value = "18BitKanKAN"

# Row 20
$ws.Cells.Item(20, 2).Value = "HotbitToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue 20 4 "0.003378"
$ws.Cells.Item(20, 5).Value = "19HotbitTokenHTB"

# Row 21
$ws.Cells.Item(21, 2).Value = "NitroEx"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue 21 4 "0.0001501"
$ws.Cells.Item(21, 5).Value = "20NitroExNTX"

# Row 22
$ws.Cells.Item(22, 2).Value = "LEO"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue 22 4 "3.692"
$ws.Cells.Item(22, 5).Value = "21LEOLEO"

# Row 23
$ws.Cells.Item(23, 2).Value = "BTSEToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue 23 4 "2.297"
$ws.Cells.Item(23, 5).Value = "22BTSETokenBTSE"

# Row 24
$ws.Cells.Item(24, 2).Value = "One"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 24 4 "0.0006138"
$ws.Cells.Item(24, 5).Value = "23OneONEWorstin24h"

# Row 25
Set-TextValue 25 4 "0.3377"

# Row 40
Set-TextValue 40 4 "0.04609"

# Row 41
Set-TextValue 41 4 "0.007167"

# Row 42
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 42 4 "0.003904"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"

# Row 43
$ws.Cells.Item(43, 2).Value = "BKEXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 43 4 "0.1119"
$ws.Cells.Item(43, 5).Value = "42BKEXTokenBKK"

# Row 44
Set-TextValue 44 4 "0.01018"

# Row 45
Set-TextValue 45 4 "0.002974"

# Row 46
Set-TextValue 46 4 "0.00006017"

# Row 48
Set-TextValue 48 4 "0.7508"

# Row 49
Set-TextValue 49 4 "0.05501"

# Row 50
Set-TextValue 50 4 "0.00001502"
